$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.974.42"
$ws.Range("E2").Value = "  +1.03%  "

# Row 3
$ws.Range("D3").Value = "2.243.48"
$ws.Range("E3").Value = "  +2.51%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$origStyle_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "271.28"
$ws.Range("D5").Style = $origStyle_D5
$ws.Range("E5").Value = "  +4.34%  "

# Row 6
$origStyle_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.72"
$ws.Range("D6").Style = $origStyle_D6
$ws.Range("E6").Value = "  +15.94%  "

# Row 7
$ws.Range("E7").Value = "  +1.45%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$origStyle_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.639"
$ws.Range("D9").Style = $origStyle_D9
$ws.Range("E9").Value = "  +8.17%  "

# Row 10
$origStyle_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.43"
$ws.Range("D10").Style = $origStyle_D10
$ws.Range("E10").Value = "  +7.08%  "

# Row 11
$origStyle_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0959"
$ws.Range("D11").Style = $origStyle_D11
$ws.Range("E11").Value = "  +4.75%  "

# Row 12
$origStyle_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.31"
$ws.Range("D12").Style = $origStyle_D12
$ws.Range("E12").Value = "  +19.77%  "

# Row 13
$ws.Range("E13").Value = "  +1.55%  "

# Row 14
$ws.Range("E14").Value = "  +7.83%  "

# Row 15
$ws.Range("D15").Value = "2.579.26"
$ws.Range("E15").Value = "  +2.49%  "

# Row 16
$origStyle_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.819"
$ws.Range("D16").Style = $origStyle_D16
$ws.Range("E16").Value = "  +5.62%  "

# Row 17
$ws.Range("D17").Value = "2.247.74"
$ws.Range("E17").Value = "  +2.88%  "

# Row 18
$ws.Range("D18").Value = "43.928.89"
$ws.Range("E18").Value = "  +1.15%  "

# Row 19
$ws.Range("E19").Value = "  +2.53%  "

# Row 20
$origStyle_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.17"
$ws.Range("D20").Style = $origStyle_D20
$ws.Range("E20").Value = "  +4.67%  "

# Row 21
$origStyle_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.85"
$ws.Range("D21").Style = $origStyle_D21
$ws.Range("E21").Value = "  +1.63%  "

# Row 22
$origStyle_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.34"
$ws.Range("D22").Style = $origStyle_D22
$ws.Range("E22").Value = "  -3.51%  "

# Row 23
$origStyle_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.91"
$ws.Range("D23").Style = $origStyle_D23

# Row 24
$origStyle_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.15"
$ws.Range("D24").Style = $origStyle_D24
$ws.Range("E24").Value = "  +4.37%  "

# Row 26
$origStyle_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("D26").Style = $origStyle_D26
$ws.Range("E26").Value = "  +7.71%  "

# Row 27
$origStyle_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.51"
$ws.Range("D27").Style = $origStyle_D27

# Row 28
$origStyle_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.58"
$ws.Range("D28").Style = $origStyle_D28
$ws.Range("E28").Value = "  +6.21%  "

# Row 29
$origStyle_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.49"
$ws.Range("D29").Style = $origStyle_D29
$ws.Range("E29").Value = "  -3.62%  "

# Row 30
$ws.Range("E30").Value = "  +0.48%  "

# Row 31
$origStyle_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.67"
$ws.Range("D31").Style = $origStyle_D31
$ws.Range("E31").Value = "  -0.42%  "

# Row 32
$origStyle_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0910"
$ws.Range("D32").Style = $origStyle_D32
$ws.Range("E32").Value = "  +5.23%  "

# Row 33
$ws.Range("E33").Value = "  +3.07%  "

# Row 34
$origStyle_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.52"
$ws.Range("D34").Style = $origStyle_D34
$ws.Range("E34").Value = "  +3.65%  "

# Row 35
$ws.Range("E35").Value = "  +2.37%  "

# Row 36
$ws.Range("E36").Value = "  -0.72%  "

# Row 37
$origStyle_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0354"
$ws.Range("D37").Style = $origStyle_D37
$ws.Range("E37").Value = "  +1.07%  "

# Row 38
$origStyle_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.32"
$ws.Range("D38").Style = $origStyle_D38
$ws.Range("E38").Value = "  -2.85%  "

# Row 39
$origStyle_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.61"
$ws.Range("D39").Style = $origStyle_D39
$ws.Range("E39").Value = "  +26.23%  "

# Row 40
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$origStyle_D40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.227"
$ws.Range("D40").Style = $origStyle_D40
$ws.Range("E40").Value = "  +14.55%  "

# Row 41
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$origStyle_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.84"
$ws.Range("D41").Style = $origStyle_D41
$ws.Range("E41").Value = "  -1.87%  "

# Row 42
$origStyle_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.17"
$ws.Range("D42").Style = $origStyle_D42
$ws.Range("E42").Value = "  +4.44%  "

# Row 43
$origStyle_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.62"
$ws.Range("D43").Style = $origStyle_D43
$ws.Range("E43").Value = "  +0.76%  "

# Row 44
$origStyle_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.44"
$ws.Range("D44").Style = $origStyle_D44
$ws.Range("E44").Value = "  +0.15%  "

# Row 45
$origStyle_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0998"
$ws.Range("D45").Style = $origStyle_D45
$ws.Range("E45").Value = "  +1.90%  "

# Row 46
$origStyle_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.33"
$ws.Range("D46").Style = $origStyle_D46
$ws.Range("E46").Value = "  +1.24%  "

# Row 47
$origStyle_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.41"
$ws.Range("D47").Style = $origStyle_D47
$ws.Range("E47").Value = "  +2.33%  "

# Row 48
$ws.Range("E48").Value = "  +4.91%  "

# Row 49
$origStyle_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.21"
$ws.Range("D49").Style = $origStyle_D49
$ws.Range("E49").Value = "  +3.54%  "

# Row 50
$origStyle_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.443"
$ws.Range("D50").Style = $origStyle_D50
$ws.Range("E50").Value = "  +1.30%  "

# Row 51
$ws.Range("D51").Value = "2.462.40"
$ws.Range("E51").Value = "  +2.46%  "
